$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "congestion_patterns" block and the "erfStruct" block (rows 5-15)
# are reordered so that congestion_patterns now comes first, followed by
# erfStruct. Also one extra blank separator row (of the 3 blank rows
# between the erfStruct/congestion_patterns block and the knobs block)
# is removed, shifting every row below up by 2.

# Step 1: remove 2 of the 3 blank rows separating row 15 from row 19
# (knobs.link_ids), leaving a single blank row (old row 16) and shifting
# everything below up by two rows.
$ws.Rows("17:18").Delete()

# Step 2: swap the erfStruct block (rows 5:6) with the congestion_patterns
# block (rows 7:15) using a scratch area far below the used range as a
# staging ground, so the two blocks trade places without clobbering each
# other.
$ws.Range("A7:E15").Cut($ws.Range("A100"))
$ws.Range("A5:E6").Cut($ws.Range("A109"))
$ws.Range("A100:E108").Cut($ws.Range("A5"))
$ws.Range("A109:E110").Cut($ws.Range("A14"))
